$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Plan")

# Row 13: Menu Screen Music status -> Done
$ws.Range("D13").Value = "Done"

# Row 25: Enemy Movement - estimated date set, assigned to Zacari
$ws.Range("B25").Value = "12/13/2021"
$ws.Range("C25").Value = "Zacari"

# Row 39: Key Placement assigned to -> Michael
$ws.Range("C39").Value = "Michael"

# Sheet view adjustment: selection moves to F36 (and the prior scrolled
# topLeftCell pin is cleared as a natural side effect of selecting here)
[void]$ws.Range("F36").Select()
